# Insert a new data row at row 49. This shifts the existing rows 49-127
# down to 50-128 (Excel's default Insert behaviour also carries the
# formatting of the row above down into the newly inserted row), which
# reproduces every row-shift seen in the diff. We only need to populate
# the freshly inserted row 49 with its new values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("49:49").Insert()

$ws.Cells.Item(49, 1).Value2 = 6
$ws.Cells.Item(49, 2).Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(49, 3).Value2 = "Metropolitana"
$ws.Cells.Item(49, 4).Value2 = 44546
$ws.Cells.Item(49, 5).Value2 = 13
$ws.Cells.Item(49, 6).Value2 = "Fruta"
$ws.Cells.Item(49, 7).Value2 = 100101
$ws.Cells.Item(49, 8).Value2 = "Berries"
$ws.Cells.Item(49, 9).Value2 = 100101004
$ws.Cells.Item(49, 10).Value2 = "Frambuesa"
$ws.Cells.Item(49, 11).Value2 = "Sin especificar"
$ws.Cells.Item(49, 12).Value2 = "Primera"
$ws.Cells.Item(49, 13).Value2 = 350
$ws.Cells.Item(49, 14).Value2 = 7000
$ws.Cells.Item(49, 15).Value2 = 7000
$ws.Cells.Item(49, 16).Value2 = 7000
$ws.Cells.Item(49, 17).Value2 = "$/bandeja 2 kilos"
$ws.Cells.Item(49, 18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(49, 19).Value2 = 3500
$ws.Cells.Item(49, 20).Value2 = 2
